$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

Set-TextValue $ws "D2" "75.858.05"
Set-TextValue $ws "E2" "  +0.34%  "
Set-TextValue $ws "D3" "2.877.14"
Set-TextValue $ws "E3" "  +5.33%  "
Set-TextValue $ws "E4" "  -0.08%  "
Set-TextValue $ws "D5" "194.90"
Set-TextValue $ws "E5" "  +4.05%  "
Set-TextValue $ws "D6" "597.01"
Set-TextValue $ws "E6" "  +0.68%  "
Set-TextValue $ws "E7" "  -0.06%  "
Set-TextValue $ws "E8" "  +1.89%  "
Set-TextValue $ws "E9" "  -2.66%  "
Set-TextValue $ws "D10" "2.872.64"
Set-TextValue $ws "E10" "  +5.18%  "
Set-TextValue $ws "E11" "  +9.19%  "
Set-TextValue $ws "E12" "  -1.53%  "
Set-TextValue $ws "D13" "4.88"
Set-TextValue $ws "E13" "  +1.79%  "
Set-TextValue $ws "D14" "3.401.04"
Set-TextValue $ws "E14" "  +6.84%  "
Set-TextValue $ws "D15" "75.713.51"
Set-TextValue $ws "E15" "  +0.27%  "
Set-TextValue $ws "E16" "  -1.37%  "
Set-TextValue $ws "D17" "27.25"
Set-TextValue $ws "E17" "  +0.45%  "
Set-TextValue $ws "D18" "2.873.96"
Set-TextValue $ws "E18" "  +5.73%  "
Set-TextValue $ws "D19" "8.85"
Set-TextValue $ws "E19" "  -7.12%  "
Set-TextValue $ws "D20" "12.51"
Set-TextValue $ws "E20" "  +2.00%  "
Set-TextValue $ws "D21" "375.97"
Set-TextValue $ws "E21" "  -0.58%  "
Set-TextValue $ws "D22" "2.28"
Set-TextValue $ws "E22" "  -1.15%  "
Set-TextValue $ws "E23" "  +0.47%  "
Set-TextValue $ws "D24" "71.39"
Set-TextValue $ws "E24" "  +0.14%  "
Set-TextValue $ws "D25" "0.998"
Set-TextValue $ws "E25" "  -0.22%  "
Set-TextValue $ws "D26" "3.028.73"
Set-TextValue $ws "E26" "  +6.34%  "
Set-TextValue $ws "E27" "  -1.12%  "
Set-TextValue $ws "D28" "9.72"
Set-TextValue $ws "E28" "  +0.68%  "
Set-TextValue $ws "E29" "  +7.55%  "
Set-TextValue $ws "D30" "0.998"
Set-TextValue $ws "E30" "  +0.05%  "
Set-TextValue $ws "E31" "  -1.55%  "
Set-TextValue $ws "D32" "506.08"
Set-TextValue $ws "E32" "  -3.53%  "
Set-TextValue $ws "D33" "7.74"
Set-TextValue $ws "E33" "  -2.16%  "
Set-TextValue $ws "D34" "1.80"
Set-TextValue $ws "E34" "  +0.03%  "
Set-TextValue $ws "D35" "0.998"
Set-TextValue $ws "E35" "  -0.11%  "
Set-TextValue $ws "D36" "163.24"
Set-TextValue $ws "E36" "  +1.10%  "
Set-TextValue $ws "D37" "20.13"
Set-TextValue $ws "E37" "  +2.57%  "
Set-TextValue $ws "D38" "19.69"
Set-TextValue $ws "E38" "  +1.52%  "
Set-TextValue $ws "E39" "  -6.53%  "
Set-TextValue $ws "D40" "181.65"
Set-TextValue $ws "E40" "  +4.26%  "
Set-TextValue $ws "E41" "  -0.05%  "
Set-TextValue $ws "D42" "0.343"
Set-TextValue $ws "E42" "  +2.12%  "
Set-TextValue $ws "D43" "4.97"
Set-TextValue $ws "E43" "  -2.25%  "
Set-TextValue $ws "D44" "1.67"
Set-TextValue $ws "E45" "  +6.33%  "
Set-TextValue $ws "D46" "1.21"
Set-TextValue $ws "E46" "  -1.63%  "
Set-TextValue $ws "D47" "40.43"
Set-TextValue $ws "E47" "  +3.24%  "
Set-TextValue $ws "D48" "2.34"
Set-TextValue $ws "E48" "  -3.25%  "
Set-TextValue $ws "D49" "0.577"
Set-TextValue $ws "E49" "  +4.17%  "
Set-TextValue $ws "D50" "0.666"
Set-TextValue $ws "E50" "  +11.36%  "
Set-TextValue $ws "D51" "3.75"
Set-TextValue $ws "E51" "  +0.49%  "
